$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.07"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'23.79"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.241"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05772"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'6.425"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'3.225"
$ws.Range("D7").Style = "Normal"

$ws.Range("D9").Value = "'0.8779"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.1376"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.07101"
$ws.Range("D11").Style = "Normal"

$ws.Range("D13").Value = "'0.03033"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.09328"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'3.814"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.001542"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.04710"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'0.0006026"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.006157"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'0.001258"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.004057"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'0.00008710"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'3.541"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'2.155"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.3162"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'0.1318"
$ws.Range("D26").Style = "Normal"

$ws.Range("D28").Value = "'0.0002331"
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").Value = "'0.03734"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006247"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'0.1047"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'0.002526"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.007149"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005342"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.5355"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.002468"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
